$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: TG Master table generation w/o O-/P- FA
# Remove the erroneous "T" (TG) marker placed in column F for the O-/P- fatty acid rows (57-62)
$ws.Range("F57:F62").ClearContents()

# UI update: scroll position / selection as left by the user
$ws.Application.Goto($ws.Range("A32"), $true)
$ws.Range("K50").Select()
